$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Real time (minutes)" values for rows 3 and 4
$ws.Range("D3").Value = 25
$ws.Range("D4").Value = 50

# Split the "Dokerizacija aplikacije i baze" task (row 11) into two tasks:
# row 11 becomes "Dokerizacija baze" with an updated estimate,
# and a new row 12 is added for "Dokerizacija aplikacije".
$ws.Range("B11").Value = "Dokerizacija baze"
$ws.Range("C11").Value = 15

$ws.Range("B12").Value = "Dokerizacija aplikacije"
$ws.Range("C12").Value = 90

# Update selection to match the authored workbook state
$ws.Range("D14").Select()
